$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 (pushes the blank spacer row and the
# three summary rows down by one).
$ws.Rows.Item(18).Insert()

# Fill in the new data row, mirroring the other data rows above it.
$ws.Range("A18").Value = 2014
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = 22
$ws.Range("D18").Value = 0.53125
$ws.Range("E18").Value = 0.70833333333333337

# Extend the time-spent / hours-spent formulas down into the new row.
$ws.Range("F2:F18").FormulaR1C1 = "=(RC[-1]-RC[-2])*24*60"
$ws.Range("G3:G18").FormulaR1C1 = "=RC[-1]/60"

# Match formatting used by the rest of the data rows.
$ws.Range("A18").NumberFormat = "general"
$ws.Range("D18").NumberFormat = $ws.Range("D17").NumberFormat
$ws.Range("E18").NumberFormat = $ws.Range("E17").NumberFormat
$ws.Range("F18").NumberFormat = $ws.Range("F17").NumberFormat
$ws.Range("G18").NumberFormat = $ws.Range("G17").NumberFormat

# The old A17 cell had a leftover (unused) time-format style; clear it to
# match the rest of the "year" column.
$ws.Range("A17").NumberFormat = "general"

$ws.Range("F18").Select()
